$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: Approving IT Manager / Approving IT Director / Plan Owner all changed
# from "Alex Pashkevych" to "Yaroslav Masyuk"
$ws.Range("E2").Value = "Yaroslav Masyuk"
$ws.Range("F2").Value = "Yaroslav Masyuk"
$ws.Range("G2").Value = "Yaroslav Masyuk"

# Row 3: Plan Owner changed from "David Antolovich" to "Alex Pashkevych; Yaroslav Masyuk"
$ws.Range("G3").Value = "Alex Pashkevych; Yaroslav Masyuk"
